$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1024:1025").Insert()

$ws.Range("A1024").Value = 5
$ws.Range("B1024").Value = "Macroferia Regional de Talca"
$ws.Range("C1024").Value = "Maule"
$ws.Range("D1024").Value = 45132
$ws.Range("E1024").Value = 7
$ws.Range("F1024").Value = "Fruta"
$ws.Range("G1024").Value = 100108
$ws.Range("H1024").Value = "Tropicales y subtropicales"
$ws.Range("I1024").Value = 100108006
$ws.Range("J1024").Value = "Plátano"
$ws.Range("K1024").Value = "Sin especificar"
$ws.Range("L1024").Value = "Pintón"
$ws.Range("M1024").Value = 800
$ws.Range("N1024").Value = 12000
$ws.Range("O1024").Value = 12000
$ws.Range("P1024").Value = 12000
$ws.Range("Q1024").Value = "$/caja 20 kilos"
$ws.Range("R1024").Value = "Ecuador"
$ws.Range("S1024").Value = 600
$ws.Range("T1024").Value = 20

$ws.Range("A1025").Value = 5
$ws.Range("B1025").Value = "Macroferia Regional de Talca"
$ws.Range("C1025").Value = "Maule"
$ws.Range("D1025").Value = 45132
$ws.Range("E1025").Value = 7
$ws.Range("F1025").Value = "Fruta"
$ws.Range("G1025").Value = 100108
$ws.Range("H1025").Value = "Tropicales y subtropicales"
$ws.Range("I1025").Value = 100108006
$ws.Range("J1025").Value = "Plátano"
$ws.Range("K1025").Value = "Sin especificar"
$ws.Range("L1025").Value = "Primera Pintón"
$ws.Range("M1025").Value = 520
$ws.Range("N1025").Value = 14000
$ws.Range("O1025").Value = 14000
$ws.Range("P1025").Value = 14000
$ws.Range("Q1025").Value = "$/caja 20 kilos"
$ws.Range("R1025").Value = "Ecuador"
$ws.Range("S1025").Value = 700
$ws.Range("T1025").Value = 20
